$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 332 (weekly update: new week's Apio prices),
# pushing the existing rows 332-394 down to 334-396.
$ws.Rows.Item(332).Insert()
$ws.Rows.Item(332).Insert()

# New row 332 - Apio, Americana (o), Primera - week of 2021-10-07 (serial 44476)
$ws.Range("A332").Value = 6
$ws.Range("B332").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C332").Value = "Metropolitana"
$ws.Range("D332").Value = 44476
$ws.Range("E332").Value = 13
$ws.Range("F332").Value = 100112017
$ws.Range("G332").Value = "Apio"
$ws.Range("H332").Value = "Americana (o)"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 2500
$ws.Range("K332").Value = 6000
$ws.Range("L332").Value = 7000
$ws.Range("M332").Value = 6600
$ws.Range("N332").Value = "`$/docena de matas"
$ws.Range("O332").Value = "Región de Coquimbo"
$ws.Range("P332").Value = 1100
$ws.Range("Q332").Value = 6
$ws.Range("R332").Value = "Hortaliza"

# New row 333 - Apio, Americana (o), Segunda - week of 2021-10-07 (serial 44476)
$ws.Range("A333").Value = 6
$ws.Range("B333").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C333").Value = "Metropolitana"
$ws.Range("D333").Value = 44476
$ws.Range("E333").Value = 13
$ws.Range("F333").Value = 100112017
$ws.Range("G333").Value = "Apio"
$ws.Range("H333").Value = "Americana (o)"
$ws.Range("I333").Value = "Segunda"
$ws.Range("J333").Value = 900
$ws.Range("K333").Value = 4000
$ws.Range("L333").Value = 5000
$ws.Range("M333").Value = 4667
$ws.Range("N333").Value = "`$/docena de matas"
$ws.Range("O333").Value = "Región de Coquimbo"
$ws.Range("P333").Value = 778
$ws.Range("Q333").Value = 6
$ws.Range("R333").Value = "Hortaliza"
